$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1029.9272
$ws.Range("J17").Value = 844.12964
$ws.Range("L17").Value = 2532.38892
$ws.Range("N17").Value = -2868.38892
$ws.Range("H18").Value = 15319.214
$ws.Range("I18").Value = 2826.6667
$ws.Range("J18").Value = 18726.273
$ws.Range("K18").Value = 2826.6667
$ws.Range("L18").Value = 18726.273
$ws.Range("M18").Value = -2542.6667
$ws.Range("N18").Value = -19294.273
$ws.Range("H129").Value = 900.0862
$ws.Range("J129").Value = 873.25
$ws.Range("L129").Value = 2619.75
$ws.Range("N129").Value = -12619.75
$ws.Range("H132").Value = 1082.4117
$ws.Range("I132").Value = 1093.8125
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 3281.4375
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -751.4375
$ws.Range("N132").Value = -7760
$ws.Range("H138").Value = 2962.25
$ws.Range("I138").Value = 3168.158
$ws.Range("K138").Value = 9504.474
$ws.Range("M138").Value = -4364.474

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4067.9216
$ws.Range("I32").Value = 2648.282
$ws.Range("J32").Value = 8681.75
$ws.Range("K32").Value = 2648.282
$ws.Range("L32").Value = 8681.75
$ws.Range("M32").Value = -2361.282
$ws.Range("N32").Value = -9255.75
$ws.Range("H61").Value = 3627.5
$ws.Range("I61").Value = 1798.3334
$ws.Range("K61").Value = 1798.3334
$ws.Range("M61").Value = -1586.3334
$ws.Range("H109").Value = 78590.664
$ws.Range("J109").Value = 78590.664
$ws.Range("L109").Value = 78590.664
$ws.Range("N109").Value = -81364.664
$ws.Range("H110").Value = 1307.6
$ws.Range("I110").Value = 928
$ws.Range("J110").Value = 1982.4445
$ws.Range("K110").Value = 928
$ws.Range("L110").Value = 1982.4445
$ws.Range("M110").Value = 1117
$ws.Range("N110").Value = -6072.4445
$ws.Range("H122").Value = 1563.1177
$ws.Range("I122").Value = 1511.5555
$ws.Range("K122").Value = 4534.666499999999
$ws.Range("M122").Value = -2084.666499999999
$ws.Range("H132").Value = 2916.3333
$ws.Range("I132").Value = 1875
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5625
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3095
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 3627.5
$ws.Range("I136").Value = 1798.3334
$ws.Range("K136").Value = 5395.0002
$ws.Range("M136").Value = -2845.0002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 489.85715
$ws.Range("I94").Value = 462
$ws.Range("J94").Value = 559.5
$ws.Range("K94").Value = 462
$ws.Range("L94").Value = 559.5
$ws.Range("M94").Value = -11
$ws.Range("N94").Value = -1461.5
$ws.Range("H108").Value = 34997.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 34997.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 34997.5
$ws.Range("N108").Value = -42677.5
$ws.Range("M108").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H132").Value = 2260.4348
$ws.Range("I132").Value = 1155.8823
$ws.Range("J132").Value = 5390
$ws.Range("K132").Value = 3467.6469
$ws.Range("L132").Value = 16170
$ws.Range("M132").Value = -937.6468999999997
$ws.Range("N132").Value = -21230
$ws.Range("H134").Value = 792.0476
$ws.Range("I134").Value = 664.8946999999999
$ws.Range("K134").Value = 1994.6841
$ws.Range("M134").Value = 540.3159000000001
$ws.Range("H141").Value = 27240
$ws.Range("J141").Value = 27240
$ws.Range("L141").Value = 27240
$ws.Range("N141").Value = -37600
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 654.9091
$ws.Range("J5").Value = 952.5
$ws.Range("L5").Value = 2857.5
$ws.Range("N5").Value = -3081.5
$ws.Range("H98").Value = 800
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2400
$ws.Range("N98").Value = -5396
$ws.Range("H122").Value = 786.44446
$ws.Range("J122").Value = 1993
$ws.Range("L122").Value = 17937
$ws.Range("N122").Value = -22837
$ws.Range("H131").Value = 11103.116
$ws.Range("J131").Value = 11924.359
$ws.Range("L131").Value = 35773.077
$ws.Range("N131").Value = -45853.077
$ws.Range("H134").Value = 1438.8096
$ws.Range("I134").Value = 1169.421
$ws.Range("J134").Value = 3998
$ws.Range("K134").Value = 3508.263
$ws.Range("L134").Value = 11994
$ws.Range("M134").Value = 1561.737
$ws.Range("N134").Value = -22134
$ws.Range("H135").Value = 654.9091
$ws.Range("J135").Value = 952.5
$ws.Range("L135").Value = 8572.5
$ws.Range("N135").Value = -13642.5
$ws.Range("M98").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1803.1428
$ws.Range("I113").Value = 1555.5
$ws.Range("J113").Value = 2133.3333
$ws.Range("K113").Value = 1555.5
$ws.Range("L113").Value = 2133.3333
$ws.Range("M113").Value = 614.5
$ws.Range("N113").Value = -6473.3333
$ws.Range("H126").Value = 61689.65
$ws.Range("I126").Value = 3081.1
$ws.Range("J126").Value = 145416.14
$ws.Range("K126").Value = 9243.299999999999
$ws.Range("L126").Value = 436248.42
$ws.Range("M126").Value = -6773.299999999999
$ws.Range("N126").Value = -441188.42

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3584
$ws.Range("I61").Value = 2876
$ws.Range("K61").Value = 2876
$ws.Range("M61").Value = -2674
$ws.Range("H68").Value = 4795.8
$ws.Range("I68").Value = 4795.8
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4795.8
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4046.8
$ws.Range("H71").Value = 4795.8
$ws.Range("I71").Value = 4795.8
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 23979
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -20235
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H113").Value = 3584
$ws.Range("I113").Value = 2876
$ws.Range("K113").Value = 2876
$ws.Range("M113").Value = -706
$ws.Range("H132").Value = 2035.1428
$ws.Range("H134").Value = 44214.5
$ws.Range("J134").Value = 44214.5
$ws.Range("L134").Value = 44214.5
$ws.Range("N134").Value = -54354.5
$ws.Range("H136").Value = 3324.238
$ws.Range("I136").Value = 2043.7778
$ws.Range("J136").Value = 4284.5835
$ws.Range("K136").Value = 6131.3334
$ws.Range("L136").Value = 12853.7505
$ws.Range("M136").Value = -3581.3334
$ws.Range("N136").Value = -17953.7505
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("N108").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 5000
$ws.Range("N48").Value = -6138
$ws.Range("H81").Value = 600
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("H84").Value = 600
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("H108").Value = 57666
$ws.Range("J108").Value = 57666
$ws.Range("L108").Value = 57666
$ws.Range("N108").Value = -65346
$ws.Range("H126").Value = 6162.44
$ws.Range("I126").Value = 5548.3125
$ws.Range("K126").Value = 16644.9375
$ws.Range("M126").Value = -14174.9375
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()
